$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (stored width ends up ~0.8333 wider than the ColumnWidth
# value supplied, so back that padding out to land on the exact target of 99)
$ws.Columns.Item(2).ColumnWidth = 98.16666666666667

# New terminology row: add the "Survey" definition in column B of row 13
$ws.Range("B13").Value = "A survey is the process of collecting, aggregating, and analyzing the responses from those questionnaires. "

# Update the active selection to reflect where review left off
$ws.Range("H12").Select()
